$wb = $excel.ActiveWorkbook

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3666.6667
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3666.6667
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3666.6667
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -4634.6667

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2225
$ws.Range("I113").Value = 1966.6666
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1966.6666
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1287.3334
$ws.Range("N113").Value = -9508

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1321.8462
$ws.Range("I45").Value = 1200.3636
$ws.Range("J45").Value = 1990
$ws.Range("K45").Value = 1200.3636
$ws.Range("L45").Value = 1990
$ws.Range("M45").Value = -823.3635999999999
$ws.Range("N45").Value = -2744

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2361.762
$ws.Range("I61").Value = 1993.9412
$ws.Range("J61").Value = 3925
$ws.Range("K61").Value = 1993.9412
$ws.Range("L61").Value = 3925
$ws.Range("M61").Value = -1781.9412

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2056
$ws.Range("I122").Value = 1832.25
$ws.Range("J122").Value = 2503.5
$ws.Range("K122").Value = 5496.75
$ws.Range("L122").Value = 7510.5
$ws.Range("M122").Value = -3046.75
$ws.Range("N122").Value = -12410.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2361.762
$ws.Range("I136").Value = 1993.9412
$ws.Range("J136").Value = 3925
$ws.Range("K136").Value = 5981.8236
$ws.Range("L136").Value = 11775
$ws.Range("M136").Value = -3431.8236

# BSM row 23
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 44900
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 44900
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 44900
$ws.Range("N23").Value = -45466

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 145628.42
$ws.Range("I86").Value = 3478.4
$ws.Range("J86").Value = 501003.5
$ws.Range("K86").Value = 3478.4
$ws.Range("L86").Value = 501003.5
$ws.Range("M86").Value = -2355.4
$ws.Range("N86").Value = -503249.5

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 145628.42
$ws.Range("I89").Value = 3478.4
$ws.Range("J89").Value = 501003.5
$ws.Range("K89").Value = 17392
$ws.Range("L89").Value = 2505017.5
$ws.Range("M89").Value = -11776
$ws.Range("N89").Value = -2516249.5

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 89.46666999999999
$ws.Range("I7").Value = 52.2
$ws.Range("J7").Value = 164
$ws.Range("K7").Value = 52.2
$ws.Range("L7").Value = 164
$ws.Range("M7").Value = 60.8
$ws.Range("N7").Value = -390

# CRP row 15
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 7333.3335
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 7333.3335
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 7333.3335
$ws.Range("N15").Value = -7673.3335

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1903.0566
$ws.Range("I31").Value = 1414.0541
$ws.Range("J31").Value = 3033.875
$ws.Range("K31").Value = 1414.0541
$ws.Range("L31").Value = 3033.875
$ws.Range("M31").Value = -1119.0541
$ws.Range("N31").Value = -3623.875

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1903.0566
$ws.Range("I34").Value = 1414.0541
$ws.Range("J34").Value = 3033.875
$ws.Range("K34").Value = 1414.0541
$ws.Range("L34").Value = 3033.875
$ws.Range("M34").Value = -1212.0541
$ws.Range("N34").Value = -3437.875

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 16307.7
$ws.Range("I94").Value = 1220
$ws.Range("J94").Value = 22773.857
$ws.Range("K94").Value = 1220
$ws.Range("L94").Value = 22773.857
$ws.Range("M94").Value = -769
$ws.Range("N94").Value = -23675.857

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2628.762
$ws.Range("I99").Value = 2565.7778
$ws.Range("J99").Value = 3006.6667
$ws.Range("K99").Value = 2565.7778
$ws.Range("L99").Value = 3006.6667
$ws.Range("M99").Value = -1067.7778
$ws.Range("N99").Value = -6002.6667

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2180.2942
$ws.Range("I122").Value = 1940.64
$ws.Range("J122").Value = 2846
$ws.Range("K122").Value = 5821.92
$ws.Range("L122").Value = 8538
$ws.Range("M122").Value = -3371.92
$ws.Range("N122").Value = -13438

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2628.762
$ws.Range("I126").Value = 2565.7778
$ws.Range("J126").Value = 3006.6667
$ws.Range("K126").Value = 7697.3334
$ws.Range("L126").Value = 9020.000100000001
$ws.Range("M126").Value = -5227.3334
$ws.Range("N126").Value = -13960.0001

# CUL row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2170.6667
$ws.Range("I103").Value = 604.8
$ws.Range("J103").Value = 10000
$ws.Range("K103").Value = 1814.4
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = -935.3999999999999
$ws.Range("N103").Value = -31758

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2382213.8
$ws.Range("I129").Value = 668.9
$ws.Range("J129").Value = 4547254.5
$ws.Range("K129").Value = 2006.7
$ws.Range("L129").Value = 13641763.5
$ws.Range("M129").Value = 2993.3
$ws.Range("N129").Value = -13651763.5

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2458.9443
$ws.Range("I102").Value = 2038.4546
$ws.Range("J102").Value = 3119.7144
$ws.Range("K102").Value = 2038.4546
$ws.Range("L102").Value = 3119.7144
$ws.Range("M102").Value = -416.4546
$ws.Range("N102").Value = -6363.7144

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4700
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4700
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 14100
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -19000

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2190
$ws.Range("I126").Value = 1362
$ws.Range("J126").Value = 3225
$ws.Range("K126").Value = 4086
$ws.Range("L126").Value = 9675
$ws.Range("M126").Value = -1616

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4972.7856
$ws.Range("I40").Value = 8951
$ws.Range("J40").Value = 3381.5
$ws.Range("K40").Value = 8951
$ws.Range("L40").Value = 3381.5
$ws.Range("M40").Value = -8815
$ws.Range("N40").Value = -3653.5

# LTW row 74
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 25000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 25000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -26996

# LTW row 77
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 25000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 75000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -84984

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 32148458
$ws.Range("I122").Value = 27782934
$ws.Range("J122").Value = 40006400
$ws.Range("K122").Value = 83348802
$ws.Range("L122").Value = 120019200
$ws.Range("M122").Value = -83346352

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 63370
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 63370
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 63370
$ws.Range("N133").Value = -68430

# LTW row 135
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 500429
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 500429
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 500429
$ws.Range("N135").Value = -510569

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 680.4783
$ws.Range("I107").Value = 658.3889
$ws.Range("J107").Value = 760
$ws.Range("K107").Value = 1975.1667
$ws.Range("L107").Value = 2280
$ws.Range("M107").Value = -55.16670000000022
$ws.Range("N107").Value = -6120

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12502569
$ws.Range("I122").Value = 13891140
$ws.Range("J122").Value = 8931960
$ws.Range("K122").Value = 41673420
$ws.Range("L122").Value = 26795880
$ws.Range("M122").Value = -41670970
$ws.Range("N122").Value = -26800780
